# Generate Report for Handoff
# The localization run was regenerated: a new source-doc GUID and a new
# content hash were produced, plus fresh handoff timestamps for zh-cn and
# de-de. Update the cell text (source of truth for the shared strings)
# and keep each hyperlink's visible "display" text in sync with the new
# file names, without touching the hyperlink targets themselves.

$oldGuid = "d91584f9-0fac-47f6-b347-aeec318a9357"
$newGuid = "93d94d35-2124-41d1-939d-605f7889888b"
$oldHash = "4e7bf8f1d7d80c8f605cbd8a8f5517c307d0c8ad"
$newHash = "bc214502c5950c20e88a5214c10b6b5eec104fc0"

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: just the source markdown file name ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"

# ---- zh-cn sheet: source file name, handoff xlf name, handoff datetime ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("C2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-09 23:04:54"

# ---- de-de sheet: source file name, handoff xlf name, handoff datetime ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("C2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-09 23:04:59"

# ---- Keep the hyperlinks' visible text in sync with the renamed files. ----
# NB: iterate with foreach (not .Item(i)/[i]) -- that's what correctly
# resolves each Hyperlink's identity in this host so the edit lands on the
# existing <hyperlink> entry instead of appending a duplicate one.
foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    }
}

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = "$newGuid.$newHash.zh-cn.xlf"
    }
}

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = "$newGuid.$newHash.de-de.xlf"
    }
}
